$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.828.57'
$ws.Range('E2').Value = '  +3.50%  '
$ws.Range('D3').Value = '3.188.93'
$ws.Range('E3').Value = '  +2.33%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''536.69'
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').Value = '''144.86'
$ws.Range('E6').Value = '  +4.76%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +3.80%  '
$ws.Range('D9').Value = '''7.30'
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('E10').Value = '  +5.23%  '
$ws.Range('D11').Value = '''0.429'
$ws.Range('E11').Value = '  +4.09%  '
$ws.Range('D12').Value = '3.735.77'
$ws.Range('E12').Value = '  +2.29%  '
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('E14').Value = '  +4.51%  '
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').Value = '59.820.47'
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '3.177.46'
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').Value = '''6.19'
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('D19').Value = '''13.06'
$ws.Range('E19').Value = '  +2.12%  '
$ws.Range('D20').Value = '''8.25'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').Value = '''381.15'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '''0.998'
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('E23').Value = '  +4.30%  '
$ws.Range('D24').Value = '''69.99'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = '''8.95'
$ws.Range('E25').Value = '  +17.41%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '''0.172'
$ws.Range('E26').Value = '  +3.20%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').Value = '0.0₃0905'
$ws.Range('E28').Value = '  +2.97%  '
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('D30').Value = '''22.47'
$ws.Range('E30').Value = '  +4.41%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').Value = '''6.15'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '''5.41'
$ws.Range('E32').Value = '  +5.21%  '
$ws.Range('D33').Value = '''1.20'
$ws.Range('E33').Value = '  +2.54%  '
$ws.Range('E34').Value = '  +4.72%  '
$ws.Range('D35').Value = '''156.36'
$ws.Range('E35').Value = '  -2.60%  '
$ws.Range('E36').Value = '  +3.03%  '
$ws.Range('D37').Value = '2.772.72'
$ws.Range('E37').Value = '  +8.29%  '
$ws.Range('D38').Value = '''25.66'
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('D39').Value = '''0.0712'
$ws.Range('E39').Value = '  +6.27%  '
$ws.Range('D40').Value = '''1.68'
$ws.Range('E40').Value = '  +2.81%  '
$ws.Range('D41').Value = '''4.28'
$ws.Range('E41').Value = '  +2.99%  '
$ws.Range('D42').Value = '''0.727'
$ws.Range('E42').Value = '  +4.32%  '
$ws.Range('D43').Value = '''39.47'
$ws.Range('E43').Value = '  +2.23%  '
$ws.Range('D44').Value = '''0.0289'
$ws.Range('E44').Value = '  +7.01%  '
$ws.Range('D45').Value = '3.230.39'
$ws.Range('E45').Value = '  +2.37%  '
$ws.Range('E46').Value = '  +2.55%  '
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('E48').Value = '  +5.69%  '
$ws.Range('D49').Value = '''20.57'
$ws.Range('E49').Value = '  +3.14%  '
$ws.Range('D50').Value = '''0.777'
$ws.Range('E50').Value = '  +4.23%  '
$ws.Range('E51').Value = '  +0.01%  '
